$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-18 17:19:05"
$wsZh.Range("H2").Value = "2016-03-18 17:19:47"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-18 17:19:13"
$wsDe.Range("H2").Value = "2016-03-18 17:20:00"
